# Apply the updated cryptocurrency market data (prices / 1h volume)
# and the two name/link/value swaps (rows 29-30, row 51) described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '59.636.31'
$ws.Cells.Item(2, 5).Value = '  -0.60%  '

$ws.Cells.Item(3, 4).Value = '2.345.13'
$ws.Cells.Item(3, 5).Value = '  -2.95%  '

$ws.Cells.Item(4, 5).Value = '  +0.00%  '

$ws.Cells.Item(5, 4).Value = '''557.56'
$ws.Cells.Item(5, 5).Value = '  +0.92%  '

$ws.Cells.Item(6, 4).Value = '''132.52'
$ws.Cells.Item(6, 5).Value = '  -3.34%  '

$ws.Cells.Item(7, 4).Value = '''1.00'
$ws.Cells.Item(7, 5).Value = '  +0.01%  '

$ws.Cells.Item(8, 4).Value = '''0.578'
$ws.Cells.Item(8, 5).Value = '  -2.51%  '

$ws.Cells.Item(9, 4).Value = '''0.104'
$ws.Cells.Item(9, 5).Value = '  -1.17%  '

$ws.Cells.Item(10, 4).Value = '''5.57'
$ws.Cells.Item(10, 5).Value = '  -1.86%  '

$ws.Cells.Item(11, 5).Value = '  +0.82%  '

$ws.Cells.Item(12, 4).Value = '''0.339'
$ws.Cells.Item(12, 5).Value = '  -3.97%  '

$ws.Cells.Item(13, 4).Value = '''23.92'
$ws.Cells.Item(13, 5).Value = '  -5.72%  '

$ws.Cells.Item(14, 4).Value = '2.782.84'
$ws.Cells.Item(14, 5).Value = '  -2.27%  '

$ws.Cells.Item(15, 4).Value = '59.611.25'
$ws.Cells.Item(15, 5).Value = '  -0.55%  '

$ws.Cells.Item(16, 4).Value = '''0.0000136'
$ws.Cells.Item(16, 5).Value = '  -0.99%  '

$ws.Cells.Item(17, 4).Value = '2.359.11'
$ws.Cells.Item(17, 5).Value = '  -2.42%  '

$ws.Cells.Item(18, 4).Value = '''10.95'
$ws.Cells.Item(18, 5).Value = '  -3.08%  '

$ws.Cells.Item(19, 4).Value = '''4.43'
$ws.Cells.Item(19, 5).Value = '  +0.33%  '

$ws.Cells.Item(20, 4).Value = '''318.12'
$ws.Cells.Item(20, 5).Value = '  -3.38%  '

$ws.Cells.Item(21, 4).Value = '''6.60'
$ws.Cells.Item(21, 5).Value = '  -0.93%  '

$ws.Cells.Item(22, 4).Value = '''0.999'
$ws.Cells.Item(22, 5).Value = '  -0.18%  '

$ws.Cells.Item(23, 4).Value = '''64.02'
$ws.Cells.Item(23, 5).Value = '  -2.92%  '

$ws.Cells.Item(24, 4).Value = '''0.172'
$ws.Cells.Item(24, 5).Value = '  -1.68%  '

$ws.Cells.Item(25, 5).Value = '  +0.16%  '

$ws.Cells.Item(26, 4).Value = '''8.31'
$ws.Cells.Item(26, 5).Value = '  -3.83%  '

$ws.Cells.Item(27, 4).Value = '''1.35'
$ws.Cells.Item(27, 5).Value = '  -1.71%  '

$ws.Cells.Item(28, 4).Value = '''1.80'
$ws.Cells.Item(28, 5).Value = '  +1.42%  '

$ws.Cells.Item(29, 2).Value = 'Monero'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(29, 4).Value = '''170.75'
$ws.Cells.Item(29, 5).Value = '  +0.96%  '

$ws.Cells.Item(30, 2).Value = 'PEPE'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(30, 4).Value = '0.0₃0749'
$ws.Cells.Item(30, 5).Value = '  -3.46%  '

$ws.Cells.Item(31, 4).Value = '''6.01'
$ws.Cells.Item(31, 5).Value = '  -0.68%  '

$ws.Cells.Item(32, 4).Value = '''1.09'
$ws.Cells.Item(32, 5).Value = '  +7.29%  '

$ws.Cells.Item(33, 4).Value = '''0.395'
$ws.Cells.Item(33, 5).Value = '  -2.45%  '

$ws.Cells.Item(34, 4).Value = '''17.99'
$ws.Cells.Item(34, 5).Value = '  -3.28%  '

$ws.Cells.Item(36, 4).Value = '''1.30'
$ws.Cells.Item(36, 5).Value = '  -1.09%  '

$ws.Cells.Item(38, 4).Value = '''4.07'
$ws.Cells.Item(38, 5).Value = '  -2.82%  '

$ws.Cells.Item(39, 4).Value = '''1.57'
$ws.Cells.Item(39, 5).Value = '  -2.02%  '

$ws.Cells.Item(40, 4).Value = '''315.20'
$ws.Cells.Item(40, 5).Value = '  -1.82%  '

$ws.Cells.Item(41, 4).Value = '''38.54'
$ws.Cells.Item(41, 5).Value = '  -2.45%  '

$ws.Cells.Item(42, 4).Value = '''144.42'
$ws.Cells.Item(42, 5).Value = '  +3.04%  '

$ws.Cells.Item(43, 4).Value = '''3.49'
$ws.Cells.Item(43, 5).Value = '  -4.77%  '

$ws.Cells.Item(44, 4).Value = '''0.0957'
$ws.Cells.Item(44, 5).Value = '  -1.29%  '

$ws.Cells.Item(45, 4).Value = '''19.12'
$ws.Cells.Item(45, 5).Value = '  -2.29%  '

$ws.Cells.Item(46, 4).Value = '''0.0505'
$ws.Cells.Item(46, 5).Value = '  -1.85%  '

$ws.Cells.Item(47, 4).Value = '''0.562'
$ws.Cells.Item(47, 5).Value = '  -2.95%  '

$ws.Cells.Item(48, 4).Value = '''0.0215'
$ws.Cells.Item(48, 5).Value = '  -3.26%  '

$ws.Cells.Item(49, 5).Value = '  +0.11%  '

$ws.Cells.Item(50, 5).Value = '  -0.19%  '

$ws.Cells.Item(51, 2).Value = 'BitgetToken'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'
$ws.Cells.Item(51, 4).Value = '''0.941'
$ws.Cells.Item(51, 5).Value = '  -0.24%  '
